$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (cells are plain text, not numbers, in the source data)
$updates = @{
    "D2" = "303.47"
    "E2" = "0.50%"
    "D3" = "37.16"
    "E3" = "6.82%"
    "D4" = "4.982"
    "E4" = "-3.87%"
    "D5" = "0.07832"
    "E5" = "1.03%"
    "D6" = "2.221"
    "E6" = "-2.47%"
    "D7" = "8.027"
    "E7" = "0.15%"
    "D8" = "4.041"
    "E8" = "0.84%"
    "D9" = "0.9148"
    "E9" = "-1.44%"
    "D10" = "0.09707"
    "E10" = "-3.64%"
    "D11" = "0.1886"
    "E11" = "3.97%"
    "D12" = "0.08568"
    "E12" = "0.48%"
    "D13" = "0.03546"
    "E13" = "1.87%"
    "D14" = "0.09956"
    "E14" = "0.43%"
    "D15" = "0.001497"
    "E15" = "0.93%"
    "D16" = "0.005659"
    "E16" = "-2.98%"
    "E17" = "-0.32%"
    "E19" = "0.71%"
    "E20" = "-2.07%"
    "D21" = "4.784"
    "E21" = "3.76%"
    "E22" = "-1.98%"
    "D23" = "0.04604"
    "E23" = "-0.39%"
    "D24" = "0.001233"
    "E24" = "0.58%"
    "D25" = "0.004787"
    "E25" = "8.29%"
    "E26" = "-8.12%"
    "E27" = "38.58%"
    "D39" = "0.01785"
    "E39" = "2.09%"
    "D40" = "0.04751"
    "E40" = "0.84%"
    "D41" = "0.008015"
    "E41" = "5.22%"
    "D42" = "0.1393"
    "E42" = "-0.90%"
    "D43" = "0.007693"
    "E43" = "11.79%"
    "D44" = "0.002162"
    "E44" = "-2.36%"
    "D45" = "0.009846"
    "E45" = "6.86%"
    "D46" = "0.00006120"
    "E46" = "2.90%"
    "D47" = "0.00000000751"
    "E47" = "-0.41%"
    "D48" = "7.618"
    "E48" = "178.82%"
    "E49" = "-0.77%"
    "D50" = "0.00002102"
    "E50" = "-0.41%"
    "D51" = "0.0002002"
    "E51" = "-0.41%"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    # Prefix with an apostrophe so Excel stores the new value as literal text
    # (matching the existing inline-string / text cell type), not a number or percentage.
    $range.Value = "'" + $updates[$cellRef]
    # Drop the quote-prefix formatting Excel applies when a value starts with an apostrophe,
    # so the cell keeps the workbook default (unstyled) look of the surrounding data cells.
    $range.ClearFormats()
}
